$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: shift EndUserPO from AC1 to AB1 (copy AA1 header format), set AC1 to EditOrderLines ---
$ws.Range("AA1").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("AB1").Value = "EndUserPO"
$ws.Range("AC1").Value = "EditOrderLines"

# --- Row 3 data: shift EPO123 from AC3 to AB3, set AC3 to new value ---
$ws.Range("AB3").Value = "EPO123"
$ws.Range("AC3").Value = "20,60,3"

# --- Row 3 C3: change order number and right-align (reuse fontId=6 format from F2/G2/K2, add right align) ---
$ws.Range("F2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "20-32122-11"
$ws.Range("C3").HorizontalAlignment = -4152

# --- Row 3 height ---
$ws.Rows("3").RowHeight = 18.75

# --- Row 6 (new): E6 right aligned empty cell (reuse default font, add right align) ---
$ws.Range("E6").HorizontalAlignment = -4152

# --- Row 10 (new): big bold heading + detail text ---
$ws.Range("B10").Value = "20-32122-11"
$ws.Range("B10").Font.Name = "Roboto"
$ws.Range("B10").Font.Size = 18
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").Font.Color = 3355443
$ws.Range("B10").VerticalAlignment = -4108
$ws.Range("B10").WrapText = $true
$ws.Range("D10").Value = "37-26283-11"
$ws.Rows("10").RowHeight = 46.5

# --- Column widths ---
$ws.Columns("T:T").ColumnWidth = 27.140625
$ws.Columns("U:U").ColumnWidth = 29
$ws.Columns("V:V").ColumnWidth = 9.42578125

# --- Sheet view / selection ---
$ws.Range("AA11").Select()
